$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.640.80"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "3.396.52"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'576.74"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "'142.66"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D9").Value = "'7.64"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").Value = "'0.385"
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("D12").Value = "3.974.96"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "'27.96"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").Value = "3.406.99"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "61.687.63"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "'6.12"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "'13.63"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'9.11"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").Value = "'388.53"
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("D22").Value = "'74.63"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "'7.38"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "'7.98"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "'1.40"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D33").Value = "'23.33"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "'6.93"
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'168.37"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'5.10"
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("D37").Value = "3.429.06"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "'0.0762"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").Value = "'26.81"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("D45").Value = "2.474.97"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").Value = "'22.66"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").Value = "'6.65"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "'2.04"
$ws.Range("E50").Value = "  -4.54%  "
$ws.Range("D51").Value = "'0.206"
$ws.Range("E51").Value = "  -1.11%  "
